$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (A5=3, B5=0) following the existing pattern
$ws.Range("A5").Value = 3.0
$ws.Range("B5").Value = 0.0
